$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptive text cells to reference FY 2012-2016 instead of FY 2011-2016
$ws.Range("A3").Value2 = "This table shows the grant awards and award dollars FDA made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the FDA page of this website."
$ws.Range("A7").Value2 = "Grant awards and award dollars FDA made for FY 2012-2016."

# Update the selected/active cell to I8, matching the saved selection state
$ws.Range("I8").Select()
